$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Cells.Item(2, 1).Value = "ECs"
    $ws.Cells.Item(2, 2).Value = "Ccl28"
    $ws.Cells.Item(2, 3).Value = "Ccr3"
    $ws.Cells.Item(2, 4).Value = "Inflammatory-Mac"
    $ws.Cells.Item(2, 5).Value = 1
    $ws.Cells.Item(2, 6).Value = 0.5
    $ws.Cells.Item(2, 7).Value = 0.1014825
    $ws.Cells.Item(2, 8).Value = 0.202965
    $ws.Cells.Item(2, 9).Value = 0.5602924329074438
    $ws.Cells.Item(2, 10).Value = 0.4593111407816425
    $ws.Cells.Item(2, 11).Value = 2
    $ws.Cells.Item(2, 12).Value = 0.6666666666666666
    $ws.Cells.Item(2, 13).Value = 0.1626153333333333
    $ws.Cells.Item(2, 14).Value = 0.487846
    $ws.Cells.Item(2, 15).Value = 0.1293260700537641
    $ws.Cells.Item(2, 16).Value = 0.1293260700537641
    $ws.Cells.Item(2, 17).Value = 0.016502610565
    $ws.Cells.Item(2, 18).Value = 0.09901566339000001
    $ws.Cells.Item(2, 19).Value = 0.07246041842878202
    $ws.Cells.Item(2, 20).Value = 0.05940090476920102

    # Row 3
    $ws.Cells.Item(3, 1).Value = "ECs"
    $ws.Cells.Item(3, 2).Value = "Ccl28"
    $ws.Cells.Item(3, 3).Value = "Ccr3"
    $ws.Cells.Item(3, 4).Value = "Neutrophils"
    $ws.Cells.Item(3, 5).Value = 1
    $ws.Cells.Item(3, 6).Value = 0.5
    $ws.Cells.Item(3, 7).Value = 0.1014825
    $ws.Cells.Item(3, 8).Value = 0.202965
    $ws.Cells.Item(3, 9).Value = 0.5602924329074438
    $ws.Cells.Item(3, 10).Value = 0.4593111407816425
    $ws.Cells.Item(3, 11).Value = 3
    $ws.Cells.Item(3, 12).Value = 1
    $ws.Cells.Item(3, 13).Value = 0.8767803333333334
    $ws.Cells.Item(3, 14).Value = 2.630341
    $ws.Cells.Item(3, 15).Value = 0.6972931302732585
    $ws.Cells.Item(3, 16).Value = 0.6972931302732585
    $ws.Cells.Item(3, 17).Value = 0.08897786017750001
    $ws.Cells.Item(3, 18).Value = 0.533867161065
    $ws.Cells.Item(3, 19).Value = 0.3906880644104511
    $ws.Cells.Item(3, 20).Value = 0.3202745031250128

    # Row 4
    $ws.Cells.Item(4, 1).Value = "ECs"
    $ws.Cells.Item(4, 2).Value = "Ccl28"
    $ws.Cells.Item(4, 3).Value = "Ccr3"
    $ws.Cells.Item(4, 4).Value = "Resolving-Mac"
    $ws.Cells.Item(4, 5).Value = 1
    $ws.Cells.Item(4, 6).Value = 0.5
    $ws.Cells.Item(4, 7).Value = 0.1014825
    $ws.Cells.Item(4, 8).Value = 0.202965
    $ws.Cells.Item(4, 9).Value = 0.5602924329074438
    $ws.Cells.Item(4, 10).Value = 0.4593111407816425
    $ws.Cells.Item(4, 11).Value = 3
    $ws.Cells.Item(4, 12).Value = 1
    $ws.Cells.Item(4, 13).Value = 0.21801
    $ws.Cells.Item(4, 14).Value = 0.65403
    $ws.Cells.Item(4, 15).Value = 0.1733807996729775
    $ws.Cells.Item(4, 16).Value = 0.1733807996729775
    $ws.Cells.Item(4, 17).Value = 0.022124199825
    $ws.Cells.Item(4, 18).Value = 0.13274519895
    $ws.Cells.Item(4, 19).Value = 0.09714395006821067
    $ws.Cells.Item(4, 20).Value = 0.0796357328874287

    # Row 5
    $ws.Cells.Item(5, 1).Value = "FAPs"
    $ws.Cells.Item(5, 2).Value = "Ccl28"
    $ws.Cells.Item(5, 3).Value = "Ccr3"
    $ws.Cells.Item(5, 4).Value = "Inflammatory-Mac"
    $ws.Cells.Item(5, 5).Value = 2
    $ws.Cells.Item(5, 6).Value = 0.6666666666666666
    $ws.Cells.Item(5, 7).Value = 0.07964166666666667
    $ws.Cells.Item(5, 8).Value = 0.238925
    $ws.Cells.Item(5, 9).Value = 0.4397075670925562
    $ws.Cells.Item(5, 10).Value = 0.5406888592183575
    $ws.Cells.Item(5, 11).Value = 2
    $ws.Cells.Item(5, 12).Value = 0.6666666666666666
    $ws.Cells.Item(5, 13).Value = 0.1626153333333333
    $ws.Cells.Item(5, 14).Value = 0.487846
    $ws.Cells.Item(5, 15).Value = 0.1293260700537641
    $ws.Cells.Item(5, 16).Value = 0.1293260700537641
    $ws.Cells.Item(5, 17).Value = 0.01295095617222222
    $ws.Cells.Item(5, 18).Value = 0.11655860555
    $ws.Cells.Item(5, 19).Value = 0.05686565162498212
    $ws.Cells.Item(5, 20).Value = 0.06992516528456312

    # Row 6
    $ws.Cells.Item(6, 1).Value = "FAPs"
    $ws.Cells.Item(6, 2).Value = "Ccl28"
    $ws.Cells.Item(6, 3).Value = "Ccr3"
    $ws.Cells.Item(6, 4).Value = "Neutrophils"
    $ws.Cells.Item(6, 5).Value = 2
    $ws.Cells.Item(6, 6).Value = 0.6666666666666666
    $ws.Cells.Item(6, 7).Value = 0.07964166666666667
    $ws.Cells.Item(6, 8).Value = 0.238925
    $ws.Cells.Item(6, 9).Value = 0.4397075670925562
    $ws.Cells.Item(6, 10).Value = 0.5406888592183575
    $ws.Cells.Item(6, 11).Value = 3
    $ws.Cells.Item(6, 12).Value = 1
    $ws.Cells.Item(6, 13).Value = 0.8767803333333334
    $ws.Cells.Item(6, 14).Value = 2.630341
    $ws.Cells.Item(6, 15).Value = 0.6972931302732585
    $ws.Cells.Item(6, 16).Value = 0.6972931302732585
    $ws.Cells.Item(6, 17).Value = 0.06982824704722222
    $ws.Cells.Item(6, 18).Value = 0.628454223425
    $ws.Cells.Item(6, 19).Value = 0.3066050658628073
    $ws.Cells.Item(6, 20).Value = 0.3770186271482457

    # Row 7
    $ws.Cells.Item(7, 1).Value = "FAPs"
    $ws.Cells.Item(7, 2).Value = "Ccl28"
    $ws.Cells.Item(7, 3).Value = "Ccr3"
    $ws.Cells.Item(7, 4).Value = "Resolving-Mac"
    $ws.Cells.Item(7, 5).Value = 2
    $ws.Cells.Item(7, 6).Value = 0.6666666666666666
    $ws.Cells.Item(7, 7).Value = 0.07964166666666667
    $ws.Cells.Item(7, 8).Value = 0.238925
    $ws.Cells.Item(7, 9).Value = 0.4397075670925562
    $ws.Cells.Item(7, 10).Value = 0.5406888592183575
    $ws.Cells.Item(7, 11).Value = 3
    $ws.Cells.Item(7, 12).Value = 1
    $ws.Cells.Item(7, 13).Value = 0.21801
    $ws.Cells.Item(7, 14).Value = 0.65403
    $ws.Cells.Item(7, 15).Value = 0.1733807996729775
    $ws.Cells.Item(7, 16).Value = 0.1733807996729775
    $ws.Cells.Item(7, 17).Value = 0.01736267975
    $ws.Cells.Item(7, 18).Value = 0.15626411775
    $ws.Cells.Item(7, 19).Value = 0.07623684960476679
    $ws.Cells.Item(7, 20).Value = 0.09374506678554875
